$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) so the new headers
# match the bold/centered/bordered look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
